$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update user credentials test data
$ws.Range("B2").Value = "Sidhi"
$ws.Range("D2").Value = "sidhip"
$ws.Range("A2").Value = "sidp@gmail.com"

# Update the active selection to A2
$ws.Range("A2").Select()
